# "Update countries & provincias Spain"
#
# The source COVID-19 country table was refreshed (new pull of case counts)
# and a handful of countries (Croacia, Ucrania, Camerun, Ghana, Nueva
# Caledonia) moved up in the ranking. Because the worksheet keeps one
# fixed row per table position (column A / "Pais" cell references do not
# change), the net effect on the sheet is:
#   * the "last updated" timestamp in A1 advances from 08:46 to 09:16
#   * a small set of rows gets refreshed Casos totales / Nuevos casos /
#     Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes
#     numbers (columns B-H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 09:16"

# --- Refreshed per-row numbers (columns B:H) ---------------------------
# Each entry is Row -> @{ Column = NewValue }
$rowUpdates = [ordered]@{
    15  = @{ B = 4517; C = 43;  E = 4483; G = 4;  H = 25 }
    17  = @{ E = 2608; G = 1;   H = 11 }
    23  = @{ D = 49;   E = 1606 }
    27  = @{ B = 1289; C = 53;  E = 1281 }
    54  = @{ B = 361;  C = 46;  D = 5;   E = 355; F = 6;  H = 1 }
    55  = @{ B = 357;  D = 100; E = 253; H = 4 }
    56  = @{ B = 352;  D = 4;   E = 348; F = 4;  H = 0 }
    57  = @{ B = 345;  D = 1;   E = 338; F = 33; H = 6 }
    70  = @{ B = 191;  C = 2;   D = 39;  E = 152 }
    83  = @{ F = 3 }
    94  = @{ C = 11;   D = 1;   E = 80;  F = 0;  H = 3 }
    95  = @{ C = 0;    D = 15;  E = 69;  F = 2 }
    96  = @{ B = 84;   C = 18;  D = 17;  E = 67 }
    97  = @{ B = 81;   D = 22;  E = 59 }
    98  = @{ B = 79;   D = 8;   E = 71;  H = 0 }
    102 = @{ C = 10;   D = 2;   E = 64;  F = 0 }
    103 = @{ B = 66;   C = 5;   D = 9;   E = 57;  F = 1; H = 0 }
    104 = @{ B = 62;   D = 0;   E = 61;  F = 4;  H = 1 }
    105 = @{ B = 59;   D = 17;  E = 42 }
    123 = @{ C = 0;    F = 0;   G = 0 }
    124 = @{ C = 5;    F = 1;   G = 1 }
    148 = @{ B = 9;    C = 1;   E = 9 }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
